$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Row 20: remove the custom "s=3" block formatting (row-level + cell-level) ---
$ws.Rows.Item(20).ClearFormats()

# --- Row 21: drop block formatting, but add a new "Hitmoti" label cell in column A (kept styled) ---
$ws.Rows.Item(21).ClearFormats()
$ws.Range("A21").Value = "Hitmoti"
$ws.Range("A21").Interior.ColorIndex = 0
$ws.Range("A21").Style = $ws.Range("B18").Style

# --- Row 22: same treatment, another "Hitmoti" label cell in column A ---
$ws.Rows.Item(22).ClearFormats()
$ws.Range("A22").Value = "Hitmoti"
$ws.Range("A22").Style = $ws.Range("B18").Style

# --- Row 23: drop per-cell formatting and remove the now-empty A23 cell entirely ---
$ws.Range("B23:D23").ClearFormats()
$ws.Range("A23").ClearFormats()
$ws.Range("A23").ClearContents()

# --- Row 24: drop per-cell formatting, set A24 to the new "Sakupinera" label ---
$ws.Range("B24:D24").ClearFormats()
$ws.Range("A24").ClearFormats()
$ws.Range("A24").Value = "Sakupinera"

# --- Rows 25-26: just drop the leftover per-cell formatting ---
$ws.Range("B25:D25").ClearFormats()
$ws.Range("B26:D26").ClearFormats()

# --- Update the saved selection on the Roadmap sheet ---
$ws.Range("G29").Select()
